$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recommendation feature needs an initial value seeded from the last two
# logged rows (row 7 = "a1", row 8 = "b2"). Duplicate them as rows 9-10.
$ws.Range("A7:CE7").Copy($ws.Range("A9:CE9"))
$ws.Range("A8:CE8").Copy($ws.Range("A10:CE10"))

Write-Host "Appended rows 9-10 (copied from rows 7-8) to $($ws.Name)"
